$wb = $excel.ActiveWorkbook

# The "all" sheet (tab index 2) holds the dialog/list/grid options table.
$ws = $wb.Worksheets.Item("all")

# Update the option strings held in D2 and B2 (D2 first so the shared-string
# table allocates indices in the same order Excel would when re-saving).
$ws.Range("D2").Value = '{label: "Текст внизу", editable: 2}'
$ws.Range("B2").Value = '{ label: "Текст  вверху", editable: "yes"}'

# Move the active selection from D2 to B2.
$ws.Activate()
$ws.Range("B2").Select()
